$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update D2 (target cluster) and recomputed TPM-derived metrics
$ws.Range("D2").Value = "ECs"
$ws.Range("G2").Value = 3.855689333333333
$ws.Range("I2").Value = 0.1513312545414155
$ws.Range("J2").Value = 0.1513312545414155
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.03825666666666667
$ws.Range("N2").Value = 0.11477
$ws.Range("O2").Value = 0.4331352014340976
$ws.Range("P2").Value = 0.4331352014340976
$ws.Range("Q2").Value = 0.1475058215955555
$ws.Range("R2").Value = 1.32755239436
$ws.Range("S2").Value = 0.06554689341907069
$ws.Range("T2").Value = 0.0655468934190707

# Row 3: was FAPs->FAPs (MuSCs before), now FAPs again (A3 changes from MuSCs to FAPs), D3 stays MuSCs->FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("G3").Value = 3.855689333333333
$ws.Range("H3").Value = 11.567068
$ws.Range("I3").Value = 0.1513312545414155
$ws.Range("J3").Value = 0.1513312545414155
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.05006833333333333
$ws.Range("N3").Value = 0.150205
$ws.Range("O3").Value = 0.5668647985659024
$ws.Range("P3").Value = 0.5668647985659024
$ws.Range("Q3").Value = 0.1930479387711111
$ws.Range("R3").Value = 1.73743144894
$ws.Range("S3").Value = 0.08578436112234479
$ws.Range("T3").Value = 0.0857843611223448

# Row 4: becomes MuSCs row with D4=ECs instead of old Resolving-Mac row
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "ECs"
$ws.Range("G4").Value = 5.360192666666666
$ws.Range("H4").Value = 16.080578
$ws.Range("I4").Value = 0.2103812342497758
$ws.Range("J4").Value = 0.2103812342497758
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.03825666666666667
$ws.Range("N4").Value = 0.11477
$ws.Range("O4").Value = 0.4331352014340976
$ws.Range("P4").Value = 0.4331352014340976
$ws.Range("Q4").Value = 0.2050631041177778
$ws.Range("R4").Value = 1.84556793706
$ws.Range("S4").Value = 0.09112351827473072
$ws.Range("T4").Value = 0.09112351827473074

# New row 5: MuSCs -> MuSCs
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Il18"
$ws.Range("C5").Value = "Il1rapl1"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 5.360192666666666
$ws.Range("H5").Value = 16.080578
$ws.Range("I5").Value = 0.2103812342497758
$ws.Range("J5").Value = 0.2103812342497758
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.05006833333333333
$ws.Range("N5").Value = 0.150205
$ws.Range("O5").Value = 0.5668647985659024
$ws.Range("P5").Value = 0.5668647985659024
$ws.Range("Q5").Value = 0.2683759131655555
$ws.Range("R5").Value = 2.41538321849
$ws.Range("S5").Value = 0.1192577159750451
$ws.Range("T5").Value = 0.1192577159750451

# New row 6: Resolving-Mac -> ECs
$ws.Range("A6").Value = "Resolving-Mac"
$ws.Range("B6").Value = "Il18"
$ws.Range("C6").Value = "Il1rapl1"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 16.26259133333333
$ws.Range("H6").Value = 48.787774
$ws.Range("I6").Value = 0.6382875112088087
$ws.Range("J6").Value = 0.6382875112088088
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.03825666666666667
$ws.Range("N6").Value = 0.11477
$ws.Range("O6").Value = 0.4331352014340976
$ws.Range("P6").Value = 0.4331352014340976
$ws.Range("Q6").Value = 0.6221525357755555
$ws.Range("R6").Value = 5.599372821979999
$ws.Range("S6").Value = 0.2764647897402962
$ws.Range("T6").Value = 0.2764647897402963

# New row 7: Resolving-Mac -> MuSCs
$ws.Range("A7").Value = "Resolving-Mac"
$ws.Range("B7").Value = "Il18"
$ws.Range("C7").Value = "Il1rapl1"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 16.26259133333333
$ws.Range("H7").Value = 48.787774
$ws.Range("I7").Value = 0.6382875112088087
$ws.Range("J7").Value = 0.6382875112088088
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.05006833333333333
$ws.Range("N7").Value = 0.150205
$ws.Range("O7").Value = 0.5668647985659024
$ws.Range("P7").Value = 0.5668647985659024
$ws.Range("Q7").Value = 0.8142408437411109
$ws.Range("R7").Value = 7.328167593669999
$ws.Range("S7").Value = 0.3618227214685125
$ws.Range("T7").Value = 0.3618227214685126
